$p = $ppt.ActivePresentation

# -----------------------------------------------------------------------
# 1) Slide 2: merge the trailing four runs of paragraph 1 ("и да може да
#    извиква всеки обект и " + "цялата информация " + "за него" + ".")
#    into a single run, keeping the text identical.
# -----------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$shp2 = $s2.Shapes.Item(2)
$tr2 = $shp2.TextFrame.TextRange

$para1 = $tr2.Paragraphs(1, 1)
$para1Text = $para1.Text
$mergeStart = "и да може да извиква"
$offset = $para1Text.IndexOf($mergeStart)
$mergedText = "и да може да извиква всеки обект и цялата информация за него."
$absStart = $para1.Start + $offset
$sub1 = $tr2.Characters($absStart, $mergedText.Length)
$sub1.Text = $mergedText

# -----------------------------------------------------------------------
# 2) Slide 2: merge the two runs "на обект" + "(" of the "Обновяване ..."
#    bullet (paragraph 7) into a single run "на обект(".
# -----------------------------------------------------------------------
$para7 = $tr2.Paragraphs(7, 1)
$para7Text = $para7.Text
$mergeStart2 = "на обект("
$offset2 = $para7Text.IndexOf($mergeStart2)
$absStart2 = $para7.Start + $offset2
$sub2 = $tr2.Characters($absStart2, $mergeStart2.Length)
$sub2.Text = $mergeStart2

# -----------------------------------------------------------------------
# 3) Slide 8: add a new "For more information read the documentation"
#    rectangle, matching the one already present on slide 4 (same
#    position/size/text). We copy the existing shape so the run/para
#    formatting (lang, dirty, smtClean) matches exactly, then re-home it
#    on slide 8 and rename it.
# -----------------------------------------------------------------------
$s4 = $p.Slides.Item(4)
$srcShape = $s4.Shapes.Item(2)
$srcShape.Copy()

$s8 = $p.Slides.Item(8)

# Burn through shape-id allocation so the pasted shape lands on id 6 /
# "Rectangle 5", matching the ids already used by the two pre-existing
# shapes on this slide (2, 3, 5).
$dummy = $s8.Shapes.AddShape(1, 0, 0, 10, 10)
$dummy.Delete()

$pastedRange = $s8.Shapes.Paste()
$newShape = $pastedRange.Item(1)
$newShape.Name = "Rectangle 5"
